# Update the "想去人数" (want-to-go count, column F) figures across all
# sheets to the latest scrape snapshot (gh-pages output regenerated at
# commit 456a3b4). Column F holds plain numeric values (not formulas),
# so each touched cell is just re-assigned its new count.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 139
$ws.Cells.Item(4, 6).Value = 1298
$ws.Cells.Item(8, 6).Value = 951
$ws.Cells.Item(10, 6).Value = 104
$ws.Cells.Item(12, 6).Value = 457
$ws.Cells.Item(14, 6).Value = 1842
$ws.Cells.Item(15, 6).Value = 4345
$ws.Cells.Item(16, 6).Value = 1264
$ws.Cells.Item(18, 6).Value = 2753
$ws.Cells.Item(20, 6).Value = 11
$ws.Cells.Item(21, 6).Value = 1114
$ws.Cells.Item(22, 6).Value = 3766
$ws.Cells.Item(23, 6).Value = 817
$ws.Cells.Item(26, 6).Value = 1518
$ws.Cells.Item(27, 6).Value = 2479
$ws.Cells.Item(29, 6).Value = 888
$ws.Cells.Item(30, 6).Value = 181
$ws.Cells.Item(31, 6).Value = 986
$ws.Cells.Item(32, 6).Value = 253
$ws.Cells.Item(33, 6).Value = 7
$ws.Cells.Item(35, 6).Value = 55
$ws.Cells.Item(36, 6).Value = 1436
$ws.Cells.Item(37, 6).Value = 2008
$ws.Cells.Item(38, 6).Value = 955
$ws.Cells.Item(40, 6).Value = 4
$ws.Cells.Item(41, 6).Value = 526
$ws.Cells.Item(42, 6).Value = 105
$ws.Cells.Item(44, 6).Value = 611
$ws.Cells.Item(45, 6).Value = 308
$ws.Cells.Item(46, 6).Value = 138
$ws.Cells.Item(48, 6).Value = 249

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(9, 6).Value = 22
$ws.Cells.Item(12, 6).Value = 129

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 497

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 497
$ws.Cells.Item(3, 6).Value = 139
$ws.Cells.Item(4, 6).Value = 1298
$ws.Cells.Item(7, 6).Value = 951
$ws.Cells.Item(10, 6).Value = 104
$ws.Cells.Item(15, 6).Value = 457
$ws.Cells.Item(16, 6).Value = 1842
$ws.Cells.Item(17, 6).Value = 4345
$ws.Cells.Item(18, 6).Value = 1264
$ws.Cells.Item(21, 6).Value = 2753
$ws.Cells.Item(22, 6).Value = 1114
$ws.Cells.Item(23, 6).Value = 3766
$ws.Cells.Item(24, 6).Value = 817
$ws.Cells.Item(27, 6).Value = 1518
$ws.Cells.Item(28, 6).Value = 2479
$ws.Cells.Item(29, 6).Value = 22
$ws.Cells.Item(33, 6).Value = 129
$ws.Cells.Item(34, 6).Value = 888
$ws.Cells.Item(35, 6).Value = 181
$ws.Cells.Item(36, 6).Value = 986
$ws.Cells.Item(37, 6).Value = 253
$ws.Cells.Item(39, 6).Value = 1436
$ws.Cells.Item(40, 6).Value = 2008
$ws.Cells.Item(41, 6).Value = 955
$ws.Cells.Item(42, 6).Value = 526
$ws.Cells.Item(43, 6).Value = 105
$ws.Cells.Item(44, 6).Value = 611
$ws.Cells.Item(45, 6).Value = 308
$ws.Cells.Item(46, 6).Value = 138
$ws.Cells.Item(48, 6).Value = 249

